# Insert a new data row just above row 173 (serial date 44596 =
# 2022-02-04, Acelga "Primera" reading for Femacal de La Calera,
# Provincia de Quillota), pushing the existing rows 173..282 down to
# 174..283 and growing the sheet's used range to A1:R283.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 173 downward (standard "insert row" semantics).
$ws.Rows.Item(173).Insert()

# Populate the newly-opened row 173 with the new record.
$ws.Cells.Item(173, 1).Value  = 3
$ws.Cells.Item(173, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(173, 3).Value  = "Coquimbo"
$ws.Cells.Item(173, 4).Value  = 44596
$ws.Cells.Item(173, 5).Value  = 5
$ws.Cells.Item(173, 6).Value  = 100112009
$ws.Cells.Item(173, 7).Value  = "Acelga"
$ws.Cells.Item(173, 8).Value  = "Sin especificar"
$ws.Cells.Item(173, 9).Value  = "Primera"
$ws.Cells.Item(173, 10).Value = 240
$ws.Cells.Item(173, 11).Value = 2300
$ws.Cells.Item(173, 12).Value = 2500
$ws.Cells.Item(173, 13).Value = 2400
$ws.Cells.Item(173, 14).Value = "`$/docena de atados (6 kilos)"
$ws.Cells.Item(173, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(173, 16).Value = 400
$ws.Cells.Item(173, 17).Value = 6
$ws.Cells.Item(173, 18).Value = "Hortaliza"
